$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.402.36'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '1.863.59'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.71'
$ws.Range('E5').Value = '  +0.90%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4624'
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07327'
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.07'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07835'
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').Value = '1.814.97'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.400'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.561'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008977'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').Value = '27.403.09'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.138'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.57'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '2.069.06'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.938'
$ws.Range('E25').Value = '  +4.75%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '151.96'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.48'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.060'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.110'
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '116.37'
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08854'
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.134'
$ws.Range('E32').Value = '  +5.26%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7675'
$ws.Range('E33').Value = '  +4.70%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.176'
$ws.Range('E34').Value = '  +3.61%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.522'
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.725'
$ws.Range('E36').Value = '  +10.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.081'
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05255'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('E40').Value = '  +1.33%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.088'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5145'
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1647'
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.430'
$ws.Range('E44').Value = '  +2.52%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.4817'
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.42'
$ws.Range('E46').Value = '  +2.09%  '
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '103.02'
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.651'
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06222'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '65.57'
$ws.Range('E51').Value = '  +1.79%  '
